$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '70.007.61'
$ws.Range('E2').Value = '  +3.63%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.795.83'
$ws.Range('E3').Value = '  +21.81%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '616.11'
$ws.Range('E5').Value = '  +6.84%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '177.73'
$ws.Range('E6').Value = '  -1.01%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.798.80'
$ws.Range('E7').Value = '  +22.03%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.14%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.545'
$ws.Range('E9').Value = '  +5.35%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.167'
$ws.Range('E10').Value = '  +9.90%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.40'
$ws.Range('E11').Value = '  -1.74%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.501'
$ws.Range('E12').Value = '  +6.89%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '40.71'
$ws.Range('E13').Value = '  +11.13%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000258'
$ws.Range('E14').Value = '  +6.41%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.423.23'
$ws.Range('E15').Value = '  +21.59%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.796.68'
$ws.Range('E16').Value = '  +21.60%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '70.136.78'
$ws.Range('E17').Value = '  +3.94%  '
$ws.Range('E18').Value = '  +1.11%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.61'
$ws.Range('E19').Value = '  +7.98%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '520.06'
$ws.Range('E20').Value = '  +7.01%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '16.78'
$ws.Range('E21').Value = '  +2.04%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.47'
$ws.Range('E22').Value = '  +21.98%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.742'
$ws.Range('E23').Value = '  +7.49%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '88.85'
$ws.Range('E24').Value = '  +6.00%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.50'
$ws.Range('E25').Value = '  +8.84%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '13.60'
$ws.Range('E26').Value = '  +6.22%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.93'
$ws.Range('E27').Value = '  +3.43%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.0000125'
$ws.Range('E29').Value = '  +31.60%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.52'
$ws.Range('E30').Value = '  +7.77%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.86'
$ws.Range('E31').Value = '  +9.28%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '7.92'
$ws.Range('E32').Value = '  -1.88%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '32.10'
$ws.Range('E33').Value = '  +14.02%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.115'
$ws.Range('E34').Value = '  +2.03%  '
$ws.Range('E35').Value = '  -0.15%  '
$ws.Range('B36').Value = 'Filecoin'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '6.22'
$ws.Range('E36').Value = '  +11.17%  '
$ws.Range('B37').Value = 'Mantle'
$ws.Range('C37').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.05'
$ws.Range('E37').Value = '  +10.41%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.344'
$ws.Range('E38').Value = '  +6.89%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.19'
$ws.Range('E39').Value = '  +8.26%  '
$ws.Range('E40').Value = '  +6.83%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '51.55'
$ws.Range('E41').Value = '  +5.00%  '
$ws.Range('B42').Value = 'Cosmos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '8.87'
$ws.Range('E42').Value = '  +6.63%  '
$ws.Range('B43').Value = 'Arweave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '44.56'
$ws.Range('E43').Value = '  -7.91%  '
$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '427.10'
$ws.Range('E44').Value = '  +13.10%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.114.75'
$ws.Range('E45').Value = '  +11.68%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.70'
$ws.Range('E46').Value = '  +0.50%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0366'
$ws.Range('E47').Value = '  +5.40%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '27.98'
$ws.Range('E48').Value = '  +4.92%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '136.95'
$ws.Range('E49').Value = '  +0.72%  '
$ws.Range('B50').Value = 'USDe'
$ws.Range('C50').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.00'
$ws.Range('E50').Value = '  +0.02%  '
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.48'
$ws.Range('E51').Value = '  +6.02%  '
